$wb = $excel.ActiveWorkbook

# --- Add the new "metadata" sheet (fetch "data" by name afterwards, since ---
# --- inserting a sheet shifts index-based references). ----------------------
$ws = $wb.Worksheets.Add()
$ws.Name = "metadata"
$ws.Outline.SummaryRow = 1
$ws.Outline.SummaryColumn = 1

$dataSheet = $wb.Worksheets.Item("data")

# --- Update the time_taken column (F) on the "data" sheet -------------------
$newTimes = @(
  "2021-10-05 14:20:07.164069",
  "2021-10-05 14:20:07.164077",
  "2021-10-05 14:20:07.164080",
  "2021-10-05 14:20:07.164083",
  "2021-10-05 14:20:07.164086",
  "2021-10-05 14:20:07.164089",
  "2021-10-05 14:20:07.164092",
  "2021-10-05 14:20:07.164095",
  "2021-10-05 14:20:07.164098",
  "2021-10-05 14:20:07.164101",
  "2021-10-05 14:20:07.164103",
  "2021-10-05 14:20:07.164106",
  "2021-10-05 14:20:07.164109",
  "2021-10-05 14:20:07.164112"
)

for ($i = 0; $i -lt $newTimes.Length; $i++) {
  $row = $i + 2
  $dataSheet.Cells.Item($row, 6).Value = $newTimes[$i]
}

# Header row (row 1), columns B..G, bold + bordered + center/top aligned
# (matches the style already used by the header row of the "data" sheet).
$ws.Range("B1").Value = "data_name"
$ws.Range("C1").Value = "data_id"
$ws.Range("D1").Value = "data_version"
$ws.Range("E1").Value = "data_version_created"
$ws.Range("F1").Value = "panel_query_time"
$ws.Range("G1").Value = "panel_get_request"

$headerRange = $ws.Range("B1:G1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

# Data row (row 2)
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "Endocrine neoplasms"
$ws.Range("C2").Value = 648
# data_version "1.23" must stay text, not be coerced into a number.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "1.23"
$ws.Range("E2").Value = "2021-07-28T13:55:27.110878Z"
$ws.Range("F2").Value = "2021-10-05 14:20:07.160604"
$ws.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/648/?format=json"

$a2Range = $ws.Range("A2")
$a2Range.Font.Bold = $true
$a2Range.Borders.LineStyle = 1
$a2Range.HorizontalAlignment = -4108
$a2Range.VerticalAlignment = -4160

# Move "metadata" to be the sheet right after "data" (data stays sheet 1).
$ws.Move($null, $dataSheet)
